$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 216, shifting existing rows 216-330 down to 217-331
$ws.Rows.Item(216).Insert()

# Populate the new row 216 with the new data
$ws.Cells.Item(216, 1).Value = 8
$ws.Cells.Item(216, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(216, 3).Value = "Coquimbo"
$ws.Cells.Item(216, 4).Value = 45089
$ws.Cells.Item(216, 5).Value = 4
$ws.Cells.Item(216, 6).Value = 100112037
$ws.Cells.Item(216, 7).Value = "Cebollín"
$ws.Cells.Item(216, 8).Value = "Sin especificar"
$ws.Cells.Item(216, 9).Value = "Primera"
$ws.Cells.Item(216, 10).Value = 1160
$ws.Cells.Item(216, 11).Value = 1000
$ws.Cells.Item(216, 12).Value = 1200
$ws.Cells.Item(216, 13).Value = 1100
$ws.Cells.Item(216, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(216, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(216, 16).Value = 183
$ws.Cells.Item(216, 17).Value = 6
$ws.Cells.Item(216, 18).Value = "Hortaliza"
